$wb = $excel.ActiveWorkbook

# --- Text change: "Ready for handoff" -> "In Translation" -----------------
# Overview sheet: status is mirrored in columns E (zh-cn) and F (de-de)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

# zh-cn detail sheet: Status column is C
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

# de-de detail sheet: Status column is C
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- Column width change (narrower "Status" columns) -----------------------
# Original stored width 17.2159881591797 -> new stored width 13.4101845877511
# ColumnWidth = 12.5 characters lands on the nearest representable width.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
